# The workbook's text cells contain a mojibake sequence: the UTF-8 bytes of
# "±" (U+00B1, PLUS-MINUS SIGN) were re-encoded as if they were Latin-1/
# Windows-1252, producing the two-character sequence "Â±" (U+00C2 U+00B1)
# inside values like "0.823 (0.815 Â± 0.011)". This walks every sheet that
# contains the pattern (f1_score, training_time, test_time) and repairs each
# cell in place, turning "Â±" back into the single correct "±" character,
# leaving every other character/cell untouched.

$wb = $excel.ActiveWorkbook

$cAumlGrave = [char]0x00C2   # "Â"
$cPlusMinus = [char]0x00B1   # "±"
$mojibake = "$cAumlGrave$cPlusMinus"   # "Â±" (mis-decoded UTF-8 for "±")
$fixed    = "$cPlusMinus"              # "±"

$sheetNames = @("f1_score", "training_time", "test_time")
$totalFixed = 0

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.Contains($mojibake)) {
                $cell.Value = $val.Replace($mojibake, $fixed)
                $totalFixed = $totalFixed + 1
            }
        }
    }
}

"Fixed $totalFixed cells"
